$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Timings")

$data = @(
    @(3.200000000186265, 0, 0.203, 0.004999999999999977, 0.3110000000000001, 0, 0.04899999999999993, 1.291, 0.3310000001874869),
    @(3.299999999813735, 0, 0.164, 0.004000000000000004, 0.375, 0, 0.06099999999999994, 1.284, 0.3139999998893472),
    @(0.5, 0, 0.164, 0.004000000000000004, 0.375, 0, 0.06099999999999994, 1.284, 15042.77299999967),
    @(1.800000000279397, 0, 0.235, 0, 0, 0, 0, 0.155, 0.3659999997580599),
    @(0.5, 0, 0.235, 0, 0, 0, 0, 0.155, 14898.16299999984),
    @(3.700000000186265, 0, 0.277, 0.004999999999999949, 0.3770000000000001, 0, 0.06900000000000006, 1.365, 0.3320000000712753),
    @(0.5, 0, 0.277, 0.004999999999999949, 0.3770000000000001, 0, 0.06900000000000006, 1.365, 15033.67400000025),
    @(1.799999999813735, 0, 0.218, 0, 0, 0, 0, 0.176, 0.3849999998237763),
    @(0.5, 0, 0.218, 0, 0, 0, 0, 0.176, 15040.71000000022)
)

$startRow = 42
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}
